$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update D2 (FAPs -> ECs) and the numeric columns E2:T2 ---
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3304176666666667
$ws.Range("H2").Value = 0.9912529999999999
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.298807
$ws.Range("N2").Value = 0.896421
$ws.Range("O2").Value = 0.1375952161456007
$ws.Range("P2").Value = 0.1375952161456007
$ws.Range("Q2").Value = 0.09873111172366666
$ws.Range("R2").Value = 0.888580005513
$ws.Range("S2").Value = 0.1375952161456007
$ws.Range("T2").Value = 0.1375952161456007

# --- Row 3: now holds the FAPs target-cluster record ---
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pomc"
$ws.Range("C3").Value = "Mc5r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3304176666666667
$ws.Range("H3").Value = 0.9912529999999999
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.244994333333333
$ws.Range("N3").Value = 3.734983
$ws.Range("O3").Value = 0.5732973604870303
$ws.Range("P3").Value = 0.5732973604870303
$ws.Range("Q3").Value = 0.4113681226332222
$ws.Range("R3").Value = 3.702313103699
$ws.Range("S3").Value = 0.5732973604870303
$ws.Range("T3").Value = 0.5732973604870303

# --- Row 4 (new row): sCs target-cluster record ---
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pomc"
$ws.Range("C4").Value = "Mc5r"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3304176666666667
$ws.Range("H4").Value = 0.9912529999999999
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6278366666666666
$ws.Range("N4").Value = 1.88351
$ws.Range("O4").Value = 0.2891074233673691
$ws.Range("P4").Value = 0.2891074233673691
$ws.Range("Q4").Value = 0.2074483264477777
$ws.Range("R4").Value = 1.86703493803
$ws.Range("S4").Value = 0.2891074233673691
$ws.Range("T4").Value = 0.2891074233673691
